# Updated cryptos list with GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D2:D51 occasionally hold plain numeric-looking text (prices/ranks).
# Force those specific cells to Text format first so Excel keeps them as literal
# strings (matching the source data) instead of auto-converting to numbers.
$ws.Range("D4:D12").NumberFormat = "@"
$ws.Range("D14:D22").NumberFormat = "@"
$ws.Range("D24:D42").NumberFormat = "@"
$ws.Range("D44:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.571.21"
$ws.Range("E2").Value = "  -2.49%  "
$ws.Range("D3").Value = "2.005.57"
$ws.Range("E3").Value = "  -4.85%  "
$ws.Range("D4").Value = "1.015"
$ws.Range("E4").Value = "  +0.76%  "
$ws.Range("D5").Value = "331.84"
$ws.Range("D6").Value = "1.013"
$ws.Range("E6").Value = "  +0.67%  "
$ws.Range("D7").Value = "0.5029"
$ws.Range("E7").Value = "  -4.15%  "
$ws.Range("D8").Value = "0.4252"
$ws.Range("E8").Value = "  -4.26%  "
$ws.Range("D9").Value = "53.64"
$ws.Range("E9").Value = "  -2.07%  "
$ws.Range("D10").Value = "0.09159"
$ws.Range("E10").Value = "  -3.46%  "
$ws.Range("D11").Value = "1.125"
$ws.Range("E11").Value = "  -4.11%  "
$ws.Range("D12").Value = "23.59"
$ws.Range("E12").Value = "  -5.80%  "
$ws.Range("D13").Value = "2.027.33"
$ws.Range("E13").Value = "  -2.24%  "
$ws.Range("D14").Value = "8.112"
$ws.Range("E14").Value = "  -7.38%  "
$ws.Range("D15").Value = "6.546"
$ws.Range("E15").Value = "  -5.49%  "
$ws.Range("D16").Value = "95.92"
$ws.Range("E16").Value = "  -5.80%  "
$ws.Range("D17").Value = "1.015"
$ws.Range("E17").Value = "  +0.66%  "
$ws.Range("D18").Value = "0.00001125"
$ws.Range("E18").Value = "  -3.71%  "
$ws.Range("D19").Value = "0.06669"
$ws.Range("E19").Value = "  -0.87%  "
$ws.Range("D20").Value = "19.93"
$ws.Range("E20").Value = "  -6.22%  "
$ws.Range("D21").Value = "1.012"
$ws.Range("E21").Value = "  +0.49%  "
$ws.Range("D22").Value = "6.006"
$ws.Range("E22").Value = "  -4.75%  "
$ws.Range("D23").Value = "29.618.82"
$ws.Range("E23").Value = "  -2.49%  "
$ws.Range("D24").Value = "12.00"
$ws.Range("E24").Value = "  -5.10%  "
$ws.Range("D25").Value = "2.283"
$ws.Range("E25").Value = "  -1.31%  "
$ws.Range("D26").Value = "159.48"
$ws.Range("E26").Value = "  -2.71%  "
$ws.Range("D27").Value = "20.80"
$ws.Range("E27").Value = "  -5.65%  "
$ws.Range("D28").Value = "6.476"
$ws.Range("E28").Value = "  -6.40%  "
$ws.Range("D29").Value = "2.344"
$ws.Range("E29").Value = "  -7.71%  "
$ws.Range("D30").Value = "128.60"
$ws.Range("D31").Value = "1.061"
$ws.Range("E31").Value = "  -7.62%  "
$ws.Range("D32").Value = "1.589"
$ws.Range("E32").Value = "  -8.60%  "
$ws.Range("D33").Value = "0.09966"
$ws.Range("E33").Value = "  -5.63%  "
$ws.Range("D34").Value = "5.891"
$ws.Range("E34").Value = "  -5.97%  "
$ws.Range("D35").Value = "3.801"
$ws.Range("E35").Value = "  -3.15%  "
$ws.Range("D36").Value = "9.607"
$ws.Range("E36").Value = "  -8.41%  "
$ws.Range("D37").Value = "0.02473"
$ws.Range("E37").Value = "  -6.00%  "
$ws.Range("D38").Value = "1.330"
$ws.Range("E38").Value = "  -1.38%  "
$ws.Range("D39").Value = "0.06408"
$ws.Range("E39").Value = "  -5.92%  "
$ws.Range("D40").Value = "0.6592"
$ws.Range("E40").Value = "  -6.28%  "
$ws.Range("D41").Value = "11.82"
$ws.Range("E41").Value = "  -6.07%  "
$ws.Range("D42").Value = "0.2080"
$ws.Range("E42").Value = "  -6.74%  "
$ws.Range("D44").Value = "0.6379"
$ws.Range("E44").Value = "  -6.80%  "
$ws.Range("D45").Value = "13.69"
$ws.Range("E45").Value = "  -5.89%  "
$ws.Range("D46").Value = "2.214"
$ws.Range("E46").Value = "  -6.10%  "
$ws.Range("D47").Value = "1.289"
$ws.Range("E47").Value = "  -5.04%  "
$ws.Range("D48").Value = "3.542"
$ws.Range("E48").Value = "  -2.87%  "
$ws.Range("D49").Value = "0.07012"
$ws.Range("E49").Value = "  -3.10%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.00000000325"
$ws.Range("E50").Value = "  -6.10%  "
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").Value = "1.140"
$ws.Range("E51").Value = "  -5.13%  "
